# Update the "random_sd_PhenTempBirds" table to also include PdeltaAIC
# as a covariate for the CG path. The "GR<-Pop_mean", "GR<-Trait_mean"
# and "Trait_mean<-det_Clim" relation rows move up, and the
# "GR<-det_Clim" / "Ind_GR<-det_Clim" rows get refreshed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextNumber {
    param($cell, [string]$text)
    # Cell values such as "9.640689e-03" look like numbers in scientific
    # notation. Assigning them through .Value/.Value2 makes Excel parse
    # them into a real floating point number (and pick up a "numeric"
    # style in the process). Entering them as a literal-string formula
    # and then converting that formula to its static result via a
    # copy / paste-values round trip keeps them as plain text cells
    # (shared string, no cell style) - exactly like the source file.
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

$excel.CutCopyMode = 0

# Row 2: GR<-Pop_mean
$ws.Range("A2").Value = "GR<-Pop_mean"
Set-TextNumber $ws.Range("B2") "9.640689e-03"
Set-TextNumber $ws.Range("C2") "5.352278e-02"
Set-TextNumber $ws.Range("D2") "6.948554e-15"

# Row 3: GR<-Trait_mean
$ws.Range("A3").Value = "GR<-Trait_mean"
Set-TextNumber $ws.Range("B3") "1.390531e-02"
Set-TextNumber $ws.Range("C3") "7.982369e-15"
Set-TextNumber $ws.Range("D3") "2.057394e-02"

# Row 4: Trait_mean<-det_Clim
$ws.Range("A4").Value = "Trait_mean<-det_Clim"
Set-TextNumber $ws.Range("B4") "1.586625e-01"
Set-TextNumber $ws.Range("C4") "2.039911e-14"
Set-TextNumber $ws.Range("D4") "1.806780e-01"

# Row 5: GR<-det_Clim (refreshed values)
$ws.Range("A5").Value = "GR<-det_Clim"
Set-TextNumber $ws.Range("B5") "9.806033e-03"
Set-TextNumber $ws.Range("C5") "3.411436e-13"
Set-TextNumber $ws.Range("D5") "7.796021e-03"

# Row 6: Ind_GR<-det_Clim (refreshed values)
$ws.Range("A6").Value = "Ind_GR<-det_Clim"
Set-TextNumber $ws.Range("B6") "4.587966e-04"
Set-TextNumber $ws.Range("C6") "0.000000e+00"
Set-TextNumber $ws.Range("D6") "0.000000e+00"

# Row 7: Tot_GR<-det_Clim
$ws.Range("A7").Value = "Tot_GR<-det_Clim"
Set-TextNumber $ws.Range("B7") "0.000000e+00"
Set-TextNumber $ws.Range("C7") "0.000000e+00"
Set-TextNumber $ws.Range("D7") "0.000000e+00"

$excel.CutCopyMode = 0
